$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.593.63'
$ws.Range('E2').Value = '  +1.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.761.60'
$ws.Range('E3').Value = '  -1.37%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '335.76'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3835'
$ws.Range('E7').Value = '  +1.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3410'
$ws.Range('E8').Value = '  -0.38%  '
$ws.Range('E9').Value = '  -3.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.136'
$ws.Range('E10').Value = '  -5.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07369'
$ws.Range('E11').Value = '  -1.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.30'
$ws.Range('E13').Value = '  +1.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.333'
$ws.Range('E14').Value = '  -2.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.760.93'
$ws.Range('E15').Value = '  -1.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.015'
$ws.Range('E16').Value = '  -0.93%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001074'
$ws.Range('E17').Value = '  -1.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06667'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '82.06'
$ws.Range('E19').Value = '  -2.12%  '
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.34'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.383'
$ws.Range('E22').Value = '  -3.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.569.11'
$ws.Range('E23').Value = '  +1.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.06'
$ws.Range('E24').Value = '  -2.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.378'
$ws.Range('E25').Value = '  -1.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '20.65'
$ws.Range('E26').Value = '  -3.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.424'
$ws.Range('E27').Value = '  -5.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.426'
$ws.Range('E28').Value = '  -4.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '152.79'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.31'
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.960.01'
$ws.Range('E31').Value = '  -1.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.085'
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.962'
$ws.Range('E33').Value = '  -1.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08791'
$ws.Range('E34').Value = '  +1.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.71'
$ws.Range('E35').Value = '  -4.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02397'
$ws.Range('E36').Value = '  +2.53%  '
$ws.Range('B37').Value = 'TheSandbox'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6763'
$ws.Range('E37').Value = '  -2.63%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.331'
$ws.Range('E38').Value = '  -2.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2179'
$ws.Range('E39').Value = '  -1.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06275'
$ws.Range('E40').Value = '  -0.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.249'
$ws.Range('E41').Value = '  +0.89%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.500'
$ws.Range('E42').Value = '  -9.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.269'
$ws.Range('E43').Value = '  -6.25%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.08'
$ws.Range('E45').Value = '  -2.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6233'
$ws.Range('E46').Value = '  -4.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.825'
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '131.33'
$ws.Range('E48').Value = '  +1.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.066'
$ws.Range('E49').Value = '  -3.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07341'
$ws.Range('E50').Value = '  +3.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.141'
$ws.Range('E51').Value = '  +2.17%  '
